$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 188.09091
$ws.Range("I55").Value = 207.75
$ws.Range("J55").Value = 176.85715
$ws.Range("K55").Value = 207.75
$ws.Range("L55").Value = 176.85715
$ws.Range("M55").Value = 6.25
$ws.Range("N55").Value = -604.85715
$ws.Range("H129").Value = 1724.7
$ws.Range("I129").Value = 749.5714
$ws.Range("K129").Value = 2248.7142
$ws.Range("M129").Value = 2751.2858
$ws.Range("H135").Value = 1438.0625
$ws.Range("I135").Value = 1038.6666
$ws.Range("K135").Value = 9347.999400000001
$ws.Range("M135").Value = -6812.999400000001
$ws.Range("H138").Value = 5517.6523
$ws.Range("J138").Value = 7046.294
$ws.Range("L138").Value = 21138.882
$ws.Range("N138").Value = -31418.882
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4137.3716
$ws.Range("I32").Value = 4119.0312
$ws.Range("K32").Value = 4119.0312
$ws.Range("M32").Value = -3832.0312
$ws.Range("H61").Value = 5245.5654
$ws.Range("I61").Value = 4783.8
$ws.Range("J61").Value = 8324
$ws.Range("K61").Value = 4783.8
$ws.Range("L61").Value = 8324
$ws.Range("M61").Value = -4571.8
$ws.Range("N61").Value = -8748
$ws.Range("H74").Value = 1823.7391
$ws.Range("I74").Value = 1672.3
$ws.Range("J74").Value = 2833.3333
$ws.Range("K74").Value = 1672.3
$ws.Range("L74").Value = 2833.3333
$ws.Range("M74").Value = -798.3
$ws.Range("N74").Value = -4581.3333
$ws.Range("H77").Value = 1823.7391
$ws.Range("I77").Value = 1672.3
$ws.Range("J77").Value = 2833.3333
$ws.Range("K77").Value = 8361.5
$ws.Range("L77").Value = 14166.6665
$ws.Range("M77").Value = -3993.5
$ws.Range("N77").Value = -22902.6665
$ws.Range("H114").Value = 30398
$ws.Range("J114").Value = 30398
$ws.Range("L114").Value = 30398
$ws.Range("N114").Value = -39076
$ws.Range("H118").Value = 199874.75
$ws.Range("J118").Value = 199874.75
$ws.Range("L118").Value = 199874.75
$ws.Range("N118").Value = -203188.75
$ws.Range("H136").Value = 5245.5654
$ws.Range("I136").Value = 4783.8
$ws.Range("J136").Value = 8324
$ws.Range("K136").Value = 14351.4
$ws.Range("L136").Value = 24972
$ws.Range("M136").Value = -11801.4
$ws.Range("N136").Value = -30072
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3596.6667
$ws.Range("I20").Value = 3500
$ws.Range("J20").Value = 3645
$ws.Range("K20").Value = 3500
$ws.Range("L20").Value = 3645
$ws.Range("M20").Value = -3253
$ws.Range("N20").Value = -4139
$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1862
$ws.Range("I16").Value = 804.4
$ws.Range("K16").Value = 804.4
$ws.Range("M16").Value = -517.4
$ws.Range("H31").Value = 5122.7617
$ws.Range("I31").Value = 2319.2307
$ws.Range("K31").Value = 2319.2307
$ws.Range("M31").Value = -2024.2307
$ws.Range("H34").Value = 5122.7617
$ws.Range("I34").Value = 2319.2307
$ws.Range("K34").Value = 2319.2307
$ws.Range("M34").Value = -2117.2307
$ws.Range("H58").Value = 2918.25
$ws.Range("I58").Value = 3024.6667
$ws.Range("J58").Value = 2599
$ws.Range("K58").Value = 3024.6667
$ws.Range("L58").Value = 2599
$ws.Range("M58").Value = -2821.6667
$ws.Range("N58").Value = -3005
$ws.Range("H62").Value = 6907
$ws.Range("I62").Value = 9210
$ws.Range("J62").Value = 3068.6667
$ws.Range("K62").Value = 9210
$ws.Range("L62").Value = 3068.6667
$ws.Range("M62").Value = -8586
$ws.Range("N62").Value = -4316.6667
$ws.Range("H65").Value = 6907
$ws.Range("I65").Value = 9210
$ws.Range("J65").Value = 3068.6667
$ws.Range("K65").Value = 46050
$ws.Range("L65").Value = 15343.3335
$ws.Range("M65").Value = -42930
$ws.Range("N65").Value = -21583.3335
$ws.Range("H113").Value = 1862
$ws.Range("I113").Value = 804.4
$ws.Range("K113").Value = 804.4
$ws.Range("M113").Value = 1365.6
$ws.Range("H132").Value = 2506
$ws.Range("I132").Value = 2529.0625
$ws.Range("K132").Value = 7587.1875
$ws.Range("M132").Value = -5057.1875
$ws.Range("H136").Value = 2918.25
$ws.Range("I136").Value = 3024.6667
$ws.Range("J136").Value = 2599
$ws.Range("K136").Value = 9074.000100000001
$ws.Range("L136").Value = 7797
$ws.Range("M136").Value = -6524.000100000001
$ws.Range("N136").Value = -12897
$ws.Range("H140").Value = 34500
$ws.Range("I140").Value = 34500
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 34500
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -29320
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 123665.11
$ws.Range("J141").Value = 123665.11
$ws.Range("L141").Value = 123665.11
$ws.Range("N141").Value = -134025.11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 344109.2
$ws.Range("I4").Value = 136846.11
$ws.Range("J4").Value = 10500000
$ws.Range("K4").Value = 410538.33
$ws.Range("L4").Value = 31500000
$ws.Range("M4").Value = -410426.33
$ws.Range("N4").Value = -31500224
$ws.Range("H55").Value = 833811.8
$ws.Range("I55").Value = 833811.8
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 2501435.4
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2501258.4
$ws.Range("N55").ClearContents()
$ws.Range("H92").Value = 4116.5
$ws.Range("J92").Value = 3139.8
$ws.Range("L92").Value = 9419.400000000001
$ws.Range("N92").Value = -11915.4
$ws.Range("H113").Value = 1203.6428
$ws.Range("I113").Value = 405.77777
$ws.Range("J113").Value = 2639.8
$ws.Range("K113").Value = 1217.33331
$ws.Range("L113").Value = 7919.400000000001
$ws.Range("M113").Value = 952.66669
$ws.Range("N113").Value = -12259.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 703688.7
$ws.Range("I7").Value = 1168481.4
$ws.Range("J7").Value = 6499.75
$ws.Range("K7").Value = 1168481.4
$ws.Range("L7").Value = 6499.75
$ws.Range("M7").Value = -1168369.4
$ws.Range("N7").Value = -6723.75
$ws.Range("H8").Value = 703688.7
$ws.Range("I8").Value = 1168481.4
$ws.Range("J8").Value = 6499.75
$ws.Range("K8").Value = 1168481.4
$ws.Range("L8").Value = 6499.75
$ws.Range("M8").Value = -1168342.4
$ws.Range("N8").Value = -6777.75
$ws.Range("H29").Value = 5000253.5
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 500
$ws.Range("N29").Value = -1080
$ws.Range("H70").Value = 22226656
$ws.Range("I70").Value = 33337662
$ws.Range("J70").Value = 4646.6
$ws.Range("K70").Value = 33337662
$ws.Range("L70").Value = 4646.6
$ws.Range("M70").Value = -33337392
$ws.Range("N70").Value = -5186.6
$ws.Range("H73").Value = 22226656
$ws.Range("I73").Value = 33337662
$ws.Range("J73").Value = 4646.6
$ws.Range("K73").Value = 33337662
$ws.Range("L73").Value = 4646.6
$ws.Range("M73").Value = -33336726
$ws.Range("N73").Value = -6518.6
$ws.Range("H97").Value = 297.69232
$ws.Range("I97").Value = 301.6
$ws.Range("J97").Value = 200
$ws.Range("K97").Value = 301.6
$ws.Range("L97").Value = 200
$ws.Range("M97").Value = 194.4
$ws.Range("N97").Value = -1192
$ws.Range("H102").Value = 3864.2727
$ws.Range("I102").Value = 3723
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 3723
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = -2101
$ws.Range("N102").Value = -7744
$ws.Range("H122").Value = 8333.432000000001
$ws.Range("J122").Value = 7185.25
$ws.Range("L122").Value = 21555.75
$ws.Range("N122").Value = -26455.75
$ws.Range("H132").Value = 2051.56
$ws.Range("I132").Value = 2075.7273
$ws.Range("J132").Value = 1874.3334
$ws.Range("K132").Value = 6227.1819
$ws.Range("L132").Value = 5623.0002
$ws.Range("M132").Value = -3697.1819
$ws.Range("N132").Value = -10683.0002
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1498.3334
$ws.Range("I16").Value = 1247.5
$ws.Range("K16").Value = 1247.5
$ws.Range("M16").Value = -1077.5
$ws.Range("H22").Value = 3167.7812
$ws.Range("I22").Value = 2425.7334
$ws.Range("J22").Value = 3822.5293
$ws.Range("K22").Value = 2425.7334
$ws.Range("L22").Value = 3822.5293
$ws.Range("M22").Value = -2130.7334
$ws.Range("N22").Value = -4412.5293
$ws.Range("H27").Value = 3167.7812
$ws.Range("I27").Value = 2425.7334
$ws.Range("J27").Value = 3822.5293
$ws.Range("K27").Value = 2425.7334
$ws.Range("L27").Value = 3822.5293
$ws.Range("M27").Value = -2318.7334
$ws.Range("N27").Value = -4036.5293
$ws.Range("H38").Value = 32995
$ws.Range("J38").Value = 32995
$ws.Range("L38").Value = 32995
$ws.Range("N38").Value = -33815
$ws.Range("H132").Value = 7334.6665
$ws.Range("I132").Value = 7334.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22003.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -19473.9995
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -16108
$ws.Range("H107").Value = 2621
$ws.Range("I107").Value = 4363
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 13089
$ws.Range("L107").Value = 5250
$ws.Range("M107").Value = -11169
$ws.Range("N107").Value = -9090
$ws.Range("H132").Value = 8107.048
$ws.Range("I132").Value = 8212.5
$ws.Range("K132").Value = 24637.5
$ws.Range("M132").Value = -22107.5
$ws.Range("H136").Value = 3573.6667
$ws.Range("I136").Value = 3588.4
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 10765.2
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -8215.200000000001
$ws.Range("N136").Value = -15600
